# Auto-generated Excel COM-interop script to apply numeric cell updates
# to the Louisoix_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 975.5454999999999
$ws.Range("I8").Value = 975.5454999999999
$ws.Range("K8").Value = 2926.6365
$ws.Range("M8").Value = -2787.6365
$ws.Range("H28").Value = 688.4054
$ws.Range("I28").Value = 462.06668
$ws.Range("K28").Value = 462.06668
$ws.Range("M28").Value = 22.93331999999998
$ws.Range("H33").Value = 157.35715
$ws.Range("I33").Value = 157.35715
$ws.Range("K33").Value = 157.35715
$ws.Range("M33").Value = 71.64285000000001
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H69").Value = 11249
$ws.Range("I69").Value = 7498.5
$ws.Range("J69").Value = 14999.5
$ws.Range("K69").Value = 22495.5
$ws.Range("L69").Value = 44998.5
$ws.Range("M69").Value = -21621.5
$ws.Range("N69").Value = -46746.5
$ws.Range("H72").Value = 11249
$ws.Range("I72").Value = 7498.5
$ws.Range("J72").Value = 14999.5
$ws.Range("K72").Value = 67486.5
$ws.Range("L72").Value = 134995.5
$ws.Range("M72").Value = -63118.5
$ws.Range("N72").Value = -143731.5
$ws.Range("H76").Value = 3758.2
$ws.Range("J76").Value = 3792.8
$ws.Range("L76").Value = 3792.8
$ws.Range("N76").Value = -4422.8
$ws.Range("H79").Value = 3758.2
$ws.Range("J79").Value = 3792.8
$ws.Range("L79").Value = 3792.8
$ws.Range("N79").Value = -5976.8
$ws.Range("H111").Value = 2334.353
$ws.Range("J111").Value = 2404.8
$ws.Range("L111").Value = 7214.400000000001
$ws.Range("N111").Value = -13348.4
$ws.Range("H137").Value = 5999
$ws.Range("I137").Value = 6999
$ws.Range("K137").Value = 20997
$ws.Range("M137").Value = -18447

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20598.963
$ws.Range("I32").Value = 28430
$ws.Range("K32").Value = 28430
$ws.Range("M32").Value = -28143
$ws.Range("H61").Value = 4118.25
$ws.Range("I61").Value = 3604.75
$ws.Range("K61").Value = 3604.75
$ws.Range("M61").Value = -3392.75
$ws.Range("H74").Value = 2552.6785
$ws.Range("I74").Value = 2436.4583
$ws.Range("J74").Value = 3250
$ws.Range("K74").Value = 2436.4583
$ws.Range("L74").Value = 3250
$ws.Range("M74").Value = -1562.4583
$ws.Range("N74").Value = -4998
$ws.Range("H77").Value = 2552.6785
$ws.Range("I77").Value = 2436.4583
$ws.Range("J77").Value = 3250
$ws.Range("K77").Value = 12182.2915
$ws.Range("L77").Value = 16250
$ws.Range("M77").Value = -7814.291499999999
$ws.Range("N77").Value = -24986
$ws.Range("H122").Value = 1333.0769
$ws.Range("I122").Value = 1278.9143
$ws.Range("K122").Value = 3836.7429
$ws.Range("M122").Value = -1386.7429
$ws.Range("H136").Value = 4118.25
$ws.Range("I136").Value = 3604.75
$ws.Range("K136").Value = 10814.25
$ws.Range("M136").Value = -8264.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50454
$ws.Range("H20").Value = 81714.234
$ws.Range("I20").Value = 104129.5
$ws.Range("K20").Value = 104129.5
$ws.Range("M20").Value = -103882.5
$ws.Range("H99").Value = 5917.6
$ws.Range("J99").Value = 6939.5
$ws.Range("L99").Value = 6939.5
$ws.Range("N99").Value = -9935.5
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1338.3334
$ws.Range("I31").Value = 1319.8387
$ws.Range("K31").Value = 1319.8387
$ws.Range("M31").Value = -1024.8387
$ws.Range("H34").Value = 1338.3334
$ws.Range("I34").Value = 1319.8387
$ws.Range("K34").Value = 1319.8387
$ws.Range("M34").Value = -1117.8387
$ws.Range("H105").Value = 1718.3
$ws.Range("I105").Value = 1703.6666
$ws.Range("J105").Value = 1850
$ws.Range("K105").Value = 1703.6666
$ws.Range("L105").Value = 1850
$ws.Range("M105").Value = 43.33339999999998
$ws.Range("N105").Value = -5344
$ws.Range("H107").Value = 2013.2683
$ws.Range("I107").Value = 305.3889
$ws.Range("K107").Value = 305.3889
$ws.Range("M107").Value = 1614.6111
$ws.Range("H134").Value = 255502.5
$ws.Range("I134").Value = 255502.5
$ws.Range("K134").Value = 766507.5
$ws.Range("M134").Value = -763972.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 104998.5
$ws.Range("I80").Value = 2499.5
$ws.Range("K80").Value = 7498.5
$ws.Range("M80").Value = -6562.5
$ws.Range("H83").Value = 104998.5
$ws.Range("I83").Value = 2499.5
$ws.Range("K83").Value = 22495.5
$ws.Range("M83").Value = -17815.5
$ws.Range("H132").Value = 2336.5715
$ws.Range("J132").Value = 2051.5557
$ws.Range("L132").Value = 18464.0013
$ws.Range("N132").Value = -23524.0013
$ws.Range("H141").Value = 7666.3335
$ws.Range("I141").Value = 1499.5
$ws.Range("K141").Value = 4498.5
$ws.Range("M141").Value = 681.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 26849
$ws.Range("I47").Value = 23700
$ws.Range("J47").Value = 29998
$ws.Range("K47").Value = 23700
$ws.Range("L47").Value = 29998
$ws.Range("M47").Value = -23132
$ws.Range("N47").Value = -31134
$ws.Range("H52").Value = 21437.777
$ws.Range("I52").Value = 21437.777
$ws.Range("K52").Value = 21437.777
$ws.Range("M52").Value = -21178.777
$ws.Range("H97").Value = 579.97144
$ws.Range("I97").Value = 482.0345
$ws.Range("J97").Value = 1053.3334
$ws.Range("K97").Value = 482.0345
$ws.Range("L97").Value = 1053.3334
$ws.Range("M97").Value = 13.96550000000002
$ws.Range("N97").Value = -2045.3334
$ws.Range("H122").Value = 2501.8333
$ws.Range("I122").Value = 1958.8462
$ws.Range("K122").Value = 5876.5386
$ws.Range("M122").Value = -3426.5386
$ws.Range("H132").Value = 130039.75
$ws.Range("I132").Value = 130039.75
$ws.Range("K132").Value = 390119.25
$ws.Range("M132").Value = -387589.25
$ws.Range("H141").Value = 66109.664
$ws.Range("J141").Value = 66109.664
$ws.Range("L141").Value = 66109.664
$ws.Range("N141").Value = -76469.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8171.8
$ws.Range("J46").Value = 2812.4211
$ws.Range("L46").Value = 2812.4211
$ws.Range("N46").Value = -3188.4211
$ws.Range("H48").Value = 22333
$ws.Range("I48").Value = 8500
$ws.Range("K48").Value = 8500
$ws.Range("M48").Value = -7839
$ws.Range("H61").Value = 6502
$ws.Range("I61").Value = 2952.5
$ws.Range("K61").Value = 2952.5
$ws.Range("M61").Value = -2750.5
$ws.Range("H113").Value = 6502
$ws.Range("I113").Value = 2952.5
$ws.Range("K113").Value = 2952.5
$ws.Range("M113").Value = -782.5
$ws.Range("H132").Value = 60121.715
$ws.Range("I132").Value = 94636.84
$ws.Range("K132").Value = 283910.52
$ws.Range("M132").Value = -281380.52

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3175.2222
$ws.Range("I122").Value = 3140.9375
$ws.Range("K122").Value = 9422.8125
$ws.Range("M122").Value = -6972.8125
$ws.Range("H132").Value = 58892.61
$ws.Range("I132").Value = 86414.914
$ws.Range("J132").Value = 3848
$ws.Range("K132").Value = 259244.742
$ws.Range("L132").Value = 11544
$ws.Range("M132").Value = -256714.742
$ws.Range("N132").Value = -16604
